$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Rename header cells: drop the leading capital letters so column
# names are all lower-case ("Epoch" -> "epoch", "Interval_name" -> "interval_name")
$ws.Range("A1").Value = "epoch"
$ws.Range("B1").Value = "interval_name"

$ws.Range("B2").Select()
